$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.614.24"
$ws.Range("E2").Value = "  -1.40%  "

$ws.Range("D3").Value = "2.679.28"
$ws.Range("E3").Value = "  -0.62%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.13"
$ws.Range("E5").Value = "  +0.21%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.87"
$ws.Range("E6").Value = "  +3.54%  "

$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("E8").Value = "  +0.75%  "

$ws.Range("D9").Value = "2.681.09"
$ws.Range("E9").Value = "  -0.51%  "

$ws.Range("E10").Value = "  +1.69%  "

$ws.Range("E12").Value = "  -0.19%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.23"
$ws.Range("E13").Value = "  -1.34%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.86"
$ws.Range("E14").Value = "  -1.41%  "

$ws.Range("D15").Value = "3.170.21"
$ws.Range("E15").Value = "  -0.49%  "

$ws.Range("E16").Value = "  -2.38%  "

$ws.Range("D17").Value = "67.441.46"
$ws.Range("E17").Value = "  -1.41%  "

$ws.Range("D18").Value = "2.704.85"
$ws.Range("E18").Value = "  +0.39%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.78"
$ws.Range("E19").Value = "  -0.65%  "

$ws.Range("E20").Value = "  +0.31%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "364.75"
$ws.Range("E21").Value = "  -0.76%  "

$ws.Range("E22").Value = "  -3.26%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.83"
$ws.Range("E23").Value = "  -1.12%  "

$ws.Range("E24").Value = "  -3.92%  "

$ws.Range("E25").Value = "  +0.14%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "71.18"
$ws.Range("E26").Value = "  -4.54%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.14"
$ws.Range("E27").Value = "  +1.37%  "

$ws.Range("E28").Value = "  -0.71%  "

$ws.Range("E29").Value = "  -2.33%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  -0.16%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "558.02"
$ws.Range("E31").Value = "  -2.63%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.99"
$ws.Range("E32").Value = "  -2.78%  "

$ws.Range("E33").Value = "  -3.83%  "

$ws.Range("E34").Value = "  -0.21%  "

$ws.Range("E35").Value = "  -1.78%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.06%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.55"
$ws.Range("E37").Value = "  -5.22%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.56"
$ws.Range("E38").Value = "  -1.59%  "

$ws.Range("E39").Value = "  -4.05%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.374"
$ws.Range("E40").Value = "  -1.24%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.32"
$ws.Range("E41").Value = "  -1.62%  "

$ws.Range("E42").Value = "  -4.15%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.96"
$ws.Range("E43").Value = "  +0.54%  "

$ws.Range("E44").Value = "  -5.01%  "

$ws.Range("E45").Value = "  +0.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.41"
$ws.Range("E46").Value = "  -0.44%  "

$ws.Range("E47").Value = "  -5.71%  "

$ws.Range("E48").Value = "  -1.40%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "153.30"
$ws.Range("E49").Value = "  -2.67%  "

$ws.Range("E50").Value = "  -3.27%  "

$ws.Range("E51").Value = "  -2.93%  "
